$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 86: date 2025-10-13 (serial 45943), station "四方坪站充电量(kw)"
$ws.Cells.Item(86, 1).Value = 45943
$ws.Cells.Item(86, 2).Value = "四方坪站充电量(kw)"
$row86 = @(686.52800000000013,1304.1330000000003,370.28000000000003,240.32999999999998,385.92099999999999,807.36500000000012,281.28800000000001,153.38799999999998,233.89499999999998,166.18199999999999,136.15100000000001,151.255,754.95099999999991,1017.9290000000002,512.26700000000005,185.57399999999998,269.23399999999998,205.44200000000001,175.94,29.7,92.1,76.36,113,84.878)
for ($i = 0; $i -lt $row86.Length; $i++) {
    $ws.Cells.Item(86, 3 + $i).Value = $row86[$i]
}

# New row 87: date 2025-10-13 (serial 45943), station "高岭站充电量(kw)"
$ws.Cells.Item(87, 1).Value = 45943
$ws.Cells.Item(87, 2).Value = "高岭站充电量(kw)"
$row87 = @(387.25599999999997,279.29799999999994,206.82900000000001,138.49600000000001,241.48099999999999,27.344000000000001,339.67399999999998,201.61399999999998,330.10300000000007,100.03,87.811000000000007,213.20600000000002,510.49199999999996,391.79700000000008,168.63799999999998,190.215,181.16300000000001,71.247,41.218999999999994,90.092000000000013,50.563000000000002,46.07,39.257999999999996,0)
for ($i = 0; $i -lt $row87.Length; $i++) {
    $ws.Cells.Item(87, 3 + $i).Value = $row87[$i]
}

# Update selection to match the post-edit state (activeCell B89 / sqref B89)
$ws.Range("B89").Select()
